# Draft update of KQ3 analysis
# - Fill in the previously-blank RefID ("B") column for the Perri and Park
#   studies.
# - Remove the still-in-progress Morice row at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Perri study rows (RefID 364). Some of these rows already had a styled
# (centered) blank cell in column B; others had no cell there at all, so
# make sure the center alignment style is applied before writing the value.
for ($r = 105; $r -le 141; $r++) {
    $ws.Cells.Item($r, 2).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 2).Value = 364
}

# Park study rows (RefID 286). These rows never had a column-B cell, and
# the new cells are left with the default (no) style.
for ($r = 142; $r -le 197; $r++) {
    $ws.Cells.Item($r, 2).Value = 286
}

# The draft "Morice" row at the bottom of the table is removed entirely.
$ws.Rows.Item(198).Delete()
